# Auto-generated script applying the cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws 'D2' '29.142.22'
Set-TextCell $ws 'E2' '  +0.47%  '
Set-TextCell $ws 'D3' '1.839.25'
Set-TextCell $ws 'E3' '  +0.53%  '
Set-TextCell $ws 'D4' '1.000'
Set-TextCell $ws 'E4' '  +0.13%  '
Set-TextCell $ws 'D5' '244.70'
Set-TextCell $ws 'E5' '  +1.06%  '
Set-TextCell $ws 'D6' '0.6314'
Set-TextCell $ws 'E6' '  +1.55%  '
Set-TextCell $ws 'E7' '  +0.33%  '
Set-TextCell $ws 'D8' '0.07528'
Set-TextCell $ws 'E8' '  +0.02%  '
Set-TextCell $ws 'D9' '0.2936'
Set-TextCell $ws 'E9' '  +0.94%  '
Set-TextCell $ws 'D10' '23.28'
Set-TextCell $ws 'E10' '  +3.99%  '
Set-TextCell $ws 'D11' '0.07746'
Set-TextCell $ws 'E11' '  +0.29%  '
Set-TextCell $ws 'D12' '1.837.73'
Set-TextCell $ws 'E12' '  +0.29%  '
Set-TextCell $ws 'D13' '5.007'
Set-TextCell $ws 'E13' '  +1.45%  '
Set-TextCell $ws 'D14' '0.6718'
Set-TextCell $ws 'E14' '  +1.42%  '
Set-TextCell $ws 'D15' '82.79'
Set-TextCell $ws 'E15' '  +0.54%  '
Set-TextCell $ws 'D16' '0.000009332'
Set-TextCell $ws 'E16' '  -6.53%  '
Set-TextCell $ws 'D17' '6.026'
Set-TextCell $ws 'E17' '  +0.09%  '
Set-TextCell $ws 'D18' '29.159.64'
Set-TextCell $ws 'E18' '  +0.52%  '
Set-TextCell $ws 'D19' '2.081.69'
Set-TextCell $ws 'E19' '  -0.05%  '
Set-TextCell $ws 'D20' '12.62'
Set-TextCell $ws 'E20' '  +2.67%  '
Set-TextCell $ws 'D21' '224.30'
Set-TextCell $ws 'E21' '  -0.42%  '
Set-TextCell $ws 'D22' '1.005'
Set-TextCell $ws 'E22' '  +0.55%  '
Set-TextCell $ws 'D23' '7.159'
Set-TextCell $ws 'E23' '  +0.28%  '
Set-TextCell $ws 'E24' '  +0.22%  '
Set-TextCell $ws 'D25' '159.90'
Set-TextCell $ws 'E25' '  +1.26%  '
Set-TextCell $ws 'D26' '0.1406'
Set-TextCell $ws 'E26' '  +2.89%  '
Set-TextCell $ws 'D27' '8.533'
Set-TextCell $ws 'E27' '  +1.16%  '
Set-TextCell $ws 'D28' '17.99'
Set-TextCell $ws 'E28' '  +0.56%  '
Set-TextCell $ws 'E29' '  +0.78%  '
Set-TextCell $ws 'D30' '0.05906'
Set-TextCell $ws 'E30' '  +13.84%  '
Set-TextCell $ws 'D31' '4.170'
Set-TextCell $ws 'E31' '  +2.28%  '
Set-TextCell $ws 'D32' '4.074'
Set-TextCell $ws 'E32' '  +1.48%  '
Set-TextCell $ws 'D33' '1.208'
Set-TextCell $ws 'E33' '  +1.42%  '
Set-TextCell $ws 'D34' '0.7495'
Set-TextCell $ws 'E34' '  +1.72%  '
Set-TextCell $ws 'D35' '1.855'
Set-TextCell $ws 'E35' '  +0.66%  '
Set-TextCell $ws 'D36' '1.145'
Set-TextCell $ws 'E36' '  +0.86%  '
Set-TextCell $ws 'D37' '2.678'
Set-TextCell $ws 'E37' '  -0.67%  '
Set-TextCell $ws 'D38' '1.231.64'
Set-TextCell $ws 'E38' '  -0.98%  '
Set-TextCell $ws 'D39' '2.768'
Set-TextCell $ws 'E39' '  +0.23%  '
Set-TextCell $ws 'D40' '0.01792'
Set-TextCell $ws 'E40' '  +0.61%  '
Set-TextCell $ws 'D41' '6.573'
Set-TextCell $ws 'E41' '  +4.41%  '
Set-TextCell $ws 'D42' '0.8971'
Set-TextCell $ws 'E42' '  +0.28%  '
Set-TextCell $ws 'D43' '1.004'
Set-TextCell $ws 'E43' '  +0.37%  '
Set-TextCell $ws 'D44' '102.34'
Set-TextCell $ws 'E44' '  +0.98%  '
Set-TextCell $ws 'D45' '1.980.79'
Set-TextCell $ws 'E45' '  -0.09%  '
Set-TextCell $ws 'D46' '0.07977'
Set-TextCell $ws 'E46' '  +18.45%  '
Set-TextCell $ws 'B47' 'BabyDogeCoin'
Set-TextCell $ws 'C47' 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextCell $ws 'D47' '0.00000000125'
Set-TextCell $ws 'E47' '  -2.39%  '
Set-TextCell $ws 'B48' 'Aave'
Set-TextCell $ws 'C48' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell $ws 'D48' '66.27'
Set-TextCell $ws 'E48' '  +3.87%  '
Set-TextCell $ws 'D49' '0.5097'
Set-TextCell $ws 'E49' '  -0.19%  '
Set-TextCell $ws 'D50' '0.4084'
Set-TextCell $ws 'E50' '  +1.81%  '
Set-TextCell $ws 'D51' '9.052'
Set-TextCell $ws 'E51' '  +2.43%  '
